{"js": "const replacements = [\n  [\"534\u00f78=\", \"978\u00f75=\"],\n  [\"457\u00f75=\", \"289\u00f79=\"],\n  [\"367\u00f75=\", \"518\u00f77=\"],\n  [\"934\u00f72=\", \"352\u00f78=\"],\n  [\"316\u00f76=\", \"531\u00f77=\"],\n  [\"203\u00f76=\", \"551\u00f77=\"],\n  [\"222\u00f76=\", \"588\u00f75=\"],\n  [\"664\u00f74=\", \"117\u00f77=\"],\n  [\"574\u00f77=\", \"590\u00f76=\"],\n  [\"980\u00f76=\", \"687\u00f74=\"],\n  [\"626\u00f74=\", \"894\u00f72=\"],\n  [\"410\u00f76=\", \"953\u00f74=\"],\n  [\"341\u00f77=\", \"154\u00f73=\"],\n  [\"196\u00f75=\", \"826\u00f76=\"],\n  [\"453\u00f73=\", \"187\u00f72=\"],\n  [\"616\u00f72=\", \"982\u00f78=\"],\n  [\"519\u00f78=\", \"275\u00f77=\"],\n  [\"743\u00f76=\", \"696\u00f79=\"],\n  [\"645\u00f76=\", \"186\u00f76=\"],\n  [\"279\u00f76=\", \"740\u00f72=\"],\n  [\"449\u00f79=\", \"638\u00f74=\"],\n  [\"945\u00f76=\", \"276\u00f78=\"],\n  [\"840\u00f79=\", \"829\u00f73=\"],\n  [\"727\u00f74=\", \"623\u00f72=\"],\n  [\"183\u00f77=\", \"341\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"534\u00f78=\", \"978\u00f75=\"),\n    @(\"457\u00f75=\", \"289\u00f79=\"),\n    @(\"367\u00f75=\", \"518\u00f77=\"),\n    @(\"934\u00f72=\", \"352\u00f78=\"),\n    @(\"316\u00f76=\", \"531\u00f77=\"),\n    @(\"203\u00f76=\", \"551\u00f77=\"),\n    @(\"222\u00f76=\", \"588\u00f75=\"),\n    @(\"664\u00f74=\", \"117\u00f77=\"),\n    @(\"574\u00f77=\", \"590\u00f76=\"),\n    @(\"980\u00f76=\", \"687\u00f74=\"),\n    @(\"626\u00f74=\", \"894\u00f72=\"),\n    @(\"410\u00f76=\", \"953\u00f74=\"),\n    @(\"341\u00f77=\", \"154\u00f73=\"),\n    @(\"196\u00f75=\", \"826\u00f76=\"),\n    @(\"453\u00f73=\", \"187\u00f72=\"),\n    @(\"616\u00f72=\", \"982\u00f78=\"),\n    @(\"519\u00f78=\", \"275\u00f77=\"),\n    @(\"743\u00f76=\", \"696\u00f79=\"),\n    @(\"645\u00f76=\", \"186\u00f76=\"),\n    @(\"279\u00f76=\", \"740\u00f72=\"),\n    @(\"449\u00f79=\", \"638\u00f74=\"),\n    @(\"945\u00f76=\", \"276\u00f78=\"),\n    @(\"840\u00f79=\", \"829\u00f73=\"),\n    @(\"727\u00f74=\", \"623\u00f72=\"),\n    @(\"183\u00f77=\", \"341\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
